$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1 gains four new leading cells (A1:D1); existing E1:T1 values shift so
# that every header cell now holds its own 0-based column index (A1=0 ... T1=19).
# New cells must carry the same repeating style pattern (4,4,8,4) already used
# across the row, so copy formats from existing cells before overwriting values.
# ---------------------------------------------------------------------------

# Grab style "4" (plain) and style "8" (bold, no border) from existing header
# cells, then stamp them onto the new A1:D1 cells before values are written.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> style 4
$ws.Range("B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> style 4
$ws.Range("D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> style 4

$ws.Range("G1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> style 8

$excel.CutCopyMode = 0

# Now write the header values 0..19 across A1:T1.
$headerCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
for ($i = 0; $i -lt $headerCols.Count; $i++) {
    $ws.Range($headerCols[$i] + "1").Value = $i
}

# ---------------------------------------------------------------------------
# New row 25: column totals for B:T, summing rows 3-23 of each column.
# All of row 25 uses the plain style "4" regardless of the column's own
# default style, so stamp that format across the whole row first.
# ---------------------------------------------------------------------------
$ws.Range("E1").Copy() | Out-Null
$ws.Range("B25:T25").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> style 4
$excel.CutCopyMode = 0

$ws.Range("B25").Formula = "=SUM(B3:B23)"
$ws.Range("C25:T25").Formula = "=SUM(C3:C23)"

# Drop the stale S15 selection left over from the previous save back to A1.
$ws.Range("A1").Select() | Out-Null

$wb.Save()
